$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''64.018.47'
$ws.Range("E2").Value = '  -1.06%  '

$ws.Range("D3").Value = '''3.149.22'
$ws.Range("E3").Value = '  -0.49%  '

$ws.Range("E4").Value = '  -0.01%  '

$ws.Range("D5").Value = '''600.67'
$ws.Range("E5").Value = '  -2.09%  '

$ws.Range("D6").Value = '''141.41'
$ws.Range("E6").Value = '  -3.46%  '

$ws.Range("E7").Value = '  -0.01%  '

$ws.Range("D8").Value = '''3.144.81'
$ws.Range("E8").Value = '  -0.56%  '

$ws.Range("D9").Value = '''0.527'
$ws.Range("E9").Value = '  -0.26%  '

$ws.Range("E10").Value = '  -2.58%  '

$ws.Range("D11").Value = '''5.38'
$ws.Range("E11").Value = '  -1.59%  '

$ws.Range("D12").Value = '''0.465'
$ws.Range("E12").Value = '  -1.91%  '

$ws.Range("E13").Value = '  -2.02%  '

$ws.Range("D14").Value = '''34.78'
$ws.Range("E14").Value = '  -2.91%  '

$ws.Range("D15").Value = '''3.664.91'
$ws.Range("E15").Value = '  -0.72%  '

$ws.Range("E16").Value = '  +2.65%  '

$ws.Range("D17").Value = '''64.002.71'
$ws.Range("E17").Value = '  -1.08%  '

$ws.Range("D18").Value = '''3.146.49'
$ws.Range("E18").Value = '  -0.70%  '

$ws.Range("D19").Value = '''6.79'
$ws.Range("E19").Value = '  -1.53%  '

$ws.Range("D20").Value = '''484.29'
$ws.Range("E20").Value = '  +0.89%  '

$ws.Range("D21").Value = '''14.61'
$ws.Range("E21").Value = '  -0.50%  '

$ws.Range("D22").Value = '''0.710'
$ws.Range("E22").Value = '  -1.76%  '

$ws.Range("D23").Value = '''7.72'
$ws.Range("E23").Value = '  -2.61%  '

$ws.Range("D24").Value = '''87.92'
$ws.Range("E24").Value = '  +4.34%  '

$ws.Range("D25").Value = '''13.20'
$ws.Range("E25").Value = '  -4.07%  '

$ws.Range("E26").Value = '  +0.03%  '

$ws.Range("D27").Value = '''2.76'
$ws.Range("E27").Value = '  -2.14%  '

$ws.Range("D28").Value = '''8.21'
$ws.Range("E28").Value = '  -6.21%  '

$ws.Range("E29").Value = '  -2.41%  '

$ws.Range("D30").Value = '''2.06'
$ws.Range("E30").Value = '  -2.70%  '

$ws.Range("D31").Value = '''27.36'
$ws.Range("E31").Value = '  +2.71%  '

$ws.Range("E32").Value = '  -6.75%  '

$ws.Range("E33").Value = '  -0.02%  '

$ws.Range("D34").Value = '''2.62'
$ws.Range("E34").Value = '  -2.64%  '

$ws.Range("E35").Value = '  -2.40%  '

$ws.Range("D36").Value = '''6.04'
$ws.Range("E36").Value = '  +0.44%  '

$ws.Range("D37").Value = '''52.80'
$ws.Range("E37").Value = '  -1.04%  '

$ws.Range("D38").Value = '''0.0₃0733'
$ws.Range("E38").Value = '  -7.93%  '

$ws.Range("D39").Value = '''2.91'
$ws.Range("E39").Value = '  -9.37%  '

$ws.Range("E40").Value = '  -0.63%  '

$ws.Range("D41").Value = '''431.44'
$ws.Range("E41").Value = '  -6.77%  '

$ws.Range("D42").Value = '''0.119'
$ws.Range("E42").Value = '  -0.53%  '

$ws.Range("D43").Value = '''8.34'
$ws.Range("E43").Value = '  -0.30%  '

$ws.Range("D44").Value = '''2.914.56'
$ws.Range("E44").Value = '  +1.88%  '

$ws.Range("D45").Value = '''0.260'
$ws.Range("E45").Value = '  -3.27%  '

$ws.Range("D46").Value = '''2.18'
$ws.Range("E46").Value = '  -6.39%  '

$ws.Range("D47").Value = '''2.38'
$ws.Range("E47").Value = '  -2.36%  '

$ws.Range("E48").Value = '  -0.07%  '

$ws.Range("B49").Value = 'Stellar'
$ws.Range("C49").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D49").Value = '''0.115'
$ws.Range("E49").Value = '  +0.25%  '

$ws.Range("B50").Value = 'InjectiveProtocol'
$ws.Range("C50").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D50").Value = '''25.73'
$ws.Range("E50").Value = '  -3.65%  '

$ws.Range("D51").Value = '''120.93'
$ws.Range("E51").Value = '  +0.87%  '
